# Generate Report for Archive
# - Flip the handoff status from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn/de-de status columns) and on each
#   per-locale sheet's Status column.
# - Re-run the column autosize for the (now shorter) status text, which
#   narrows those status columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status text wherever it appears.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# Resize the status columns to fit the new, shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth     = 12.5
$dede.Columns.Item(3).ColumnWidth     = 12.5
